$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "展览" (exhibitions) : A1:I12 -> A1:I13
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 117
$ws1.Range("F5").Value = 5389
$ws1.Range("F6").Value = 72
$ws1.Range("F9").Value = 2385
$ws1.Range("F12").Value = 2240

# New row 13 - copy formatting of row 12 first (A column is bold/bordered/centered)
$ws1.Range("A12:I12").Copy()
$ws1.Range("A13:I13").PasteSpecial(-4122)
$ws1.Range("A13").Value = 12
$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2024-10-03"
$ws1.Range("B13").Style = "Normal"
$ws1.Range("C13").Value = "江西·JMG（广电）第二届UP动漫游戏博览会"
$ws1.Range("D13").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws1.Range("E13").Value = "2024.10.03 09:00-10.05 18:00"
$ws1.Range("F13").Value = 32
$ws1.Range("G13").Value = 19.9
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=90599"
$ws1.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202408/2LP6dm961723428231240.jpeg"

# ------------------------------------------------------------------
# Sheet "演出" (performances)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 95

# ------------------------------------------------------------------
# Sheet "全部类型" (all types) : A1:I15 -> A1:I16
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 117
$ws4.Range("F5").Value = 5389
$ws4.Range("F6").Value = 95
$ws4.Range("F7").Value = 72
$ws4.Range("F11").Value = 2385
$ws4.Range("F15").Value = 2240

# New row 16 - copy formatting of row 15 first
$ws4.Range("A15:I15").Copy()
$ws4.Range("A16:I16").PasteSpecial(-4122)
$ws4.Range("A16").Value = 15
$ws4.Range("B16").NumberFormat = "@"
$ws4.Range("B16").Value = "2024-10-03"
$ws4.Range("B16").Style = "Normal"
$ws4.Range("C16").Value = "江西·JMG（广电）第二届UP动漫游戏博览会"
$ws4.Range("D16").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws4.Range("E16").Value = "2024.10.03 09:00-10.05 18:00"
$ws4.Range("F16").Value = 32
$ws4.Range("G16").Value = 19.9
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=90599"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202408/2LP6dm961723428231240.jpeg"
